# Insert a new weekly data row at row 356 (pushing existing rows 356-389 down
# to 357-390) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(356).Insert()

$ws.Range("A356").Value = 3
$ws.Range("B356").Value = "Femacal de La Calera"
$ws.Range("C356").Value = "Coquimbo"
$ws.Range("D356").Value = 44769
$ws.Range("E356").Value = 5
$ws.Range("F356").Value = 100112043
$ws.Range("G356").Value = "Pepino ensalada"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 125
$ws.Range("K356").Value = 17000
$ws.Range("L356").Value = 18000
$ws.Range("M356").Value = 17480
$ws.Range("N356").Value = "$/caja 70 unidades"
$ws.Range("O356").Value = "Región de Arica y Parinacota"
$ws.Range("P356").Value = 250
$ws.Range("Q356").Value = 70
$ws.Range("R356").Value = "Hortaliza"
